$d = $word.ActiveDocument

# Third table holds the test-results log; row 4 (1-indexed) is the empty
# row that needs to be filled in with the 18-AVR-2017 entry.
$t = $d.Tables.Item(3)

$cell1 = $t.Cell(4, 1)
$r1 = $cell1.Range
$r1.Text = "18-AVR-2017"
$r1.Font.Bold = $true
$r1.Font.Size = 12
$r1.Font.SizeBi = 12

$cell2 = $t.Cell(4, 2)
$r2 = $cell2.Range
$r2.Text = "Même résultat que le test précédent."
$r2.Font.Bold = $true
$r2.Font.Size = 12
$r2.Font.SizeBi = 12

$cell3 = $t.Cell(4, 3)
$r3 = $cell3.Range
$r3.Text = "non"
$r3.Font.Bold = $true
$r3.Font.Size = 12
$r3.Font.SizeBi = 12

# Move the _GoBack bookmark (previously sitting in the trailing empty
# paragraph after the table) to just after the "non" we just typed.
$cell3b = $t.Cell(4, 3)
$r3b = $cell3b.Range
$bmRange = $d.Range($r3b.End - 1, $r3b.End - 1)
$bmRange.Bookmarks.Add("_GoBack")
